$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AA147").Value = ""
$ws.Range("AA147").Font.Underline = $true
